$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for the three new BOM lines (U5 / USB1 / JST1) by
#    inserting three rows above the current blank separator row (12).
#    This pushes the separator row, the Total row and the footnote
#    row down by three (12->15, 13->16, 15->18) and keeps the SUM()
#    formula + shared-formula block range auto-adjusted.
# ------------------------------------------------------------------
$ws.Range("12:14").Insert()

# ------------------------------------------------------------------
# 2. Update reference-designator / quantity text that changed because
#    of the new components (extra resistors + capacitors added to the
#    existing R/C banks).
# ------------------------------------------------------------------
$ws.Range("B3").Value = "R1, R2, R3, R4, R5, R6, R7, R8"
$ws.Range("D3").Value = 8

$ws.Range("B4").Value = "R9, R10, R11"

$ws.Range("B5").Value = "C1, C3, C4, C5, C6"
$ws.Range("D5").Value = 5

# Part numbers picked up a manufacturer suffix.
$ws.Range("C7").Value = "MCP6G02-E/SN (SGA)"
$ws.Range("C8").Value = "MCP6S21-I/SN (PGA)"

# ------------------------------------------------------------------
# 3. Populate the three new BOM rows for the USB charger circuitry.
# ------------------------------------------------------------------
$ws.Range("B12").Value = "U5"
$ws.Range("C12").Value = "MCP73832T-2ACI/OT (battery charger)"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0.38
$ws.Range("F12").Formula = "=D12*E12"

$ws.Range("B13").Value = "USB1"
$ws.Range("C13").Value = "Micro USB socket"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.518
$ws.Range("F13").Formula = "=D13*E13"

$ws.Range("B14").Value = "JST1"
$ws.Range("C14").Value = "JST battery connector"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0.49
$ws.Range("F14").Formula = "=D14*E14"

# ------------------------------------------------------------------
# 4. Rebuild the hyperlinks collection so the newly inserted cells
#    (and the renamed op-amp part numbers) get the correct links /
#    display text while every pre-existing link keeps its target.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C11"), "http://uk.rs-online.com/web/p/microcontrollers/7421176/", "", "", "PIC24FV32KA304-I/PT")
$ws.Hyperlinks.Add($ws.Range("C8"), "http://uk.rs-online.com/web/p/programmable-gain-amplifiers/0403193/", "", "", "MCP6S21 (PGA)")
$ws.Hyperlinks.Add($ws.Range("C7"), "http://uk.rs-online.com/web/p/programmable-gain-amplifiers/0402986/", "", "", "MCP6G02 (SGA)")
$ws.Hyperlinks.Add($ws.Range("C6"), "http://uk.rs-online.com/web/p/mosfet-transistors/7527773/")
$ws.Hyperlinks.Add($ws.Range("C5"), "http://uk.rs-online.com/web/p/ceramic-multilayer-capacitors/2644371/")
$ws.Hyperlinks.Add($ws.Range("C9"), "http://uk.rs-online.com/web/p/condenser-microphone-components/7243134/")
$ws.Hyperlinks.Add($ws.Range("C3"), "http://uk.rs-online.com/web/p/surface-mount-fixed-resistors/2230477/")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://uk.rs-online.com/web/p/condenser-microphone-components/7243122/")
$ws.Hyperlinks.Add($ws.Range("C4"), "http://uk.rs-online.com/web/p/surface-mount-fixed-resistors/2230297/")
$ws.Hyperlinks.Add($ws.Range("C13"), "http://uk.rs-online.com/web/p/usb-connectors/7484121/")
$ws.Hyperlinks.Add($ws.Range("C14"), "http://uk.rs-online.com/web/p/pcb-connectors/7766274/")
$ws.Hyperlinks.Add($ws.Range("C12"), "http://uk.rs-online.com/web/p/battery-charger-ics/8061056/", "", "", "MCP73832T-2ACI/OT")

# ------------------------------------------------------------------
# 5. Cosmetics: wider B/C columns (C now best-fits its longer text),
#    and move the active-cell selection like the author left it.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 27.6
$ws.Columns.Item(3).ColumnWidth = 34.4

$ws.Range("G14").Select()
